# Generate Report for Handoff
# Updates the localization-status workbook to reflect the newly generated
# handoff report: the "8c2dbc05..." / "9841f33d..." / "dd373966..." /
# "e838e367..." files moved from "low" priority / pending handoff to "ht"
# priority with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H)
# refreshed (new handoff xliff generated for these files)
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-09-03 16:37:28"

# de-de sheet: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H)
# refreshed to the same new generation timestamp used on the Overview sheet
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-09-03 16:37:33"

# Overview sheet: rows 4-7 -> Latest HO Xliff Generate Date (G) mirrors the
# de-de handoff datetime above (same shared string in the source workbook)
$overview.Range("G4:G7").Value = "2016-09-03 16:37:33"
